$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12 (item 7): the style corrector ("correctora de estilo") delivered the
# files, so record the delivery/received dates that were previously blank.
$ws.Range("B12").Value = 42091
$ws.Range("C12").Value = 42091
$ws.Range("D12").Value = 42086

# Leave a note on B12 documenting that the style corrector delivered the
# complete files.
$comment = $ws.Range("B12").AddComment("La correctora de estilo entregó los archivos completos")
$comment.Visible = $false

# Move the active selection to G12, matching where the author left off.
$ws.Range("G12").Select()
